# Daily attendance processing - 2026-01-16 10:36:21
# Swap the order of "System" and "dnasr281@gmail.com" in the "Recorded By"
# column (G) wherever the combined value "System, dnasr281@gmail.com" is
# used, changing it to "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
